$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.443.43'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.636.36'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('E6').Value = '  +4.92%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.91'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -5.16%  '
$ws.Range('E9').Value = '  -2.23%  '
$ws.Range('E10').Value = '  -0.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0887'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.867.51'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.637.73'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.73%  '
$ws.Range('E14').Value = '  +2.75%  '
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.20'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.443.42'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.68'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0724'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.63'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.02%  '
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('E22').Value = '  -2.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.77'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.32%  '
$ws.Range('E24').Value = '  -3.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '149.68'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.58%  '
$ws.Range('E26').Value = '  -2.65%  '
$ws.Range('E27').Value = '  +1.98%  '
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.18%  '
$ws.Range('E30').Value = '  -0.67%  '
$ws.Range('E31').Value = '  -1.73%  '
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.16'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.416.34'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.24%  '
$ws.Range('E35').Value = '  +2.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.10%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.572'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0167'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.873'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.901'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +14.39%  '
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('E43').Value = '  +1.70%  '
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.98'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.777.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('E47').Value = '  -3.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '85.75'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.61%  '
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0988'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.73'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.65%  '
